$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The project row (row 2) no longer has an entry: clear the project name
# and role, and reset the counter back to 0.
$ws.Range("A2").Value = 0
$ws.Range("B2:C2").Clear()

# Move the active selection to B4 (matches the saved selection state).
$ws.Range("B4").Select()
